$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Tocantins"
$ws.Range("B2").Value = "Variação 2023/2013"
$ws.Range("C2").Value = 52.5266370526839

# Row 3
$ws.Range("A3").Value = "Mato Grosso"
$ws.Range("B3").Value = "Variação 2023/2013"
$ws.Range("C3").Value = 46.68053194226529

# Row 4
$ws.Range("A4").Value = "Roraima"
$ws.Range("B4").Value = "Variação 2023/2013"
$ws.Range("C4").Value = 36.59515826253023

# Row 5
$ws.Range("A5").Value = "Acre"
$ws.Range("B5").Value = "Variação 2023/2013"
$ws.Range("C5").Value = 34.11074656431843

# Row 6
$ws.Range("A6").Value = "Mato Grosso do Sul"
$ws.Range("B6").Value = "Variação 2023/2013"
$ws.Range("C6").Value = 30.47189773844641

# Row 7
$ws.Range("A7").Value = "Maranhão"
$ws.Range("B7").Value = "Variação 2023/2013"
$ws.Range("C7").Value = 28.8835604097077

# Row 8
$ws.Range("B8").Value = "Variação 2023/2013"
$ws.Range("C8").Value = 11.41741112685601
$ws.Range("D8").Value = "19º"

# Row 9
$ws.Range("B9").Value = "Variação 2023/2013"
$ws.Range("C9").Value = 15.67765727798931

# Row 10
$ws.Range("B10").Value = "Variação 2023/2013"
$ws.Range("C10").Value = 11.76330200396814
